$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # Overview
$ws2 = $wb.Worksheets.Item(2)  # zh-cn
$ws3 = $wb.Worksheets.Item(3)  # de-de

$oldGuid = "ed06fc48-7c75-4a54-a3b0-722d9c0141d6"
$newGuid = "24a46b52-f339-44d4-8044-bcb1c4a544e7"

$oldHashZh = "837ea8b8fb32f2a1afd1d0331692cb307a24929f"
$newHashZh = "be5090b18542a3277308b616e009b375a25223fd"

$newMdName = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHashZh.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHashZh.de-de.xlf"

$newZhDatetime = "2016-03-10 03:40:40"
$newDeDatetime = "2016-03-10 03:40:43"

# --- Sheet1 (Overview): A2 file-name cell ---
$ws1.Range("A2").Value = $newMdName
foreach ($hl in $ws1.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$A`$2") {
        $hl.TextToDisplay = $newMdName
    }
}

# --- Sheet2 (zh-cn): A2 file-name cell, C2 handoff file, D2 handoff datetime ---
$ws2.Range("A2").Value = $newMdName
$ws2.Range("C2").Value = $newZhXlfName
$ws2.Range("D2").Value = $newZhDatetime
foreach ($hl in $ws2.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq "`$C`$2") {
        $hl.TextToDisplay = $newZhXlfName
    }
}

# --- Sheet3 (de-de): A2 file-name cell, C2 handoff file, D2 handoff datetime ---
$ws3.Range("A2").Value = $newMdName
$ws3.Range("C2").Value = $newDeXlfName
$ws3.Range("D2").Value = $newDeDatetime
foreach ($hl in $ws3.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2") {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq "`$C`$2") {
        $hl.TextToDisplay = $newDeXlfName
    }
}
